$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.170.24"
$ws.Range("E2").Value = "  +0.28%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.906.31"
$ws.Range("E3").Value = "  +0.74%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.57"
$ws.Range("E5").Value = "  -0.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  +0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5236"
$ws.Range("E7").Value = "  +1.77%  "

$ws.Range("E8").Value = "  +0.38%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07251"

$ws.Range("E10").Value = "  -0.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9061"
$ws.Range("E11").Value = "  +0.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08495"
$ws.Range("E12").Value = "  +11.12%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.911.34"
$ws.Range("E13").Value = "  +1.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "96.93"
$ws.Range("E14").Value = "  +2.16%  "

$ws.Range("E15").Value = "  +0.50%  "

$ws.Range("E16").Value = "  +0.24%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008665"
$ws.Range("E17").Value = "  +2.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.56"
$ws.Range("E18").Value = "  +0.81%  "

$ws.Range("E19").Value = "  +0.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.203.73"
$ws.Range("E20").Value = "  +0.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.088"
$ws.Range("E21").Value = "  +0.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.148.99"
$ws.Range("E22").Value = "  +1.02%  "

$ws.Range("E23").Value = "  +0.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.440"
$ws.Range("E24").Value = "  +0.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.332"
$ws.Range("E25").Value = "  +2.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.99"
$ws.Range("E26").Value = "  +0.89%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.26"
$ws.Range("E27").Value = "  +1.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.746"
$ws.Range("E28").Value = "  -1.13%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.22"
$ws.Range("E29").Value = "  +0.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.929"
$ws.Range("E30").Value = "  -0.58%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.824"
$ws.Range("E31").Value = "  -0.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09305"
$ws.Range("E32").Value = "  +1.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8041"
$ws.Range("E33").Value = "  +2.63%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05062"
$ws.Range("E34").Value = "  -0.49%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.248"
$ws.Range("E35").Value = "  +0.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.445"
$ws.Range("E36").Value = "  +4.81%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.952"
$ws.Range("E37").Value = "  -1.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.612"
$ws.Range("E38").Value = "  -0.71%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5727"
$ws.Range("E39").Value = "  +2.45%  "

$ws.Range("E40").Value = "  +0.12%  "

$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("E42").Value = "  +0.31%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.641"
$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "116.05"
$ws.Range("E44").Value = "  -1.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1518"
$ws.Range("E45").Value = "  +0.53%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4868"
$ws.Range("E46").Value = "  +1.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.21"
$ws.Range("E47").Value = "  +0.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9996"
$ws.Range("E48").Value = "  +0.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.623"
$ws.Range("E49").Value = "  +1.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.69"
$ws.Range("E50").Value = "  +0.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.31"
$ws.Range("E51").Value = "  +0.45%  "
